$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---- Header row (row 1): turn the old "data-as-header" row into real
# column headers and extend it with the common property columns (G:M),
# matching the pattern used by the other property-type sheets (stock,
# other securities, insurance, debt, ...).
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

$headerCols = @("G","H","I","J","K","L","M")
$headerVals = @("property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $col = $headerCols[$i]
    $cell = $col + "1"
    $ws.Range("B1").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).Value = $headerVals[$i]
}

# ---- Data rows 2-5: the B:F values (bank name, deposit type, currency,
# owner, total) are already correct and unchanged; we only need to add
# the new G:M columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) copying format
# from the existing data columns.
$dataRows = @(2, 3, 4, 5)
$indexVals = @{ 2 = 45; 3 = 46; 4 = 47; 5 = 48 }

foreach ($r in $dataRows) {
    $srcCell = "B" + $r
    $cols = @("G","H","I","J","K","L","M")
    foreach ($col in $cols) {
        $cell = $col + $r
        $ws.Range($srcCell).Copy()
        $ws.Range($cell).PasteSpecial(-4122)
    }
    $ws.Range("G" + $r).Value = "deposit"
    $ws.Range("H" + $r).Value = "normal"
    $ws.Range("I" + $r).Value = "2011-12-06"
    $ws.Range("J" + $r).Value = "王廷升"
    $ws.Range("K" + $r).Value = 1727
    $ws.Range("L" + $r).Value = "tmp44311"
    $ws.Range("M" + $r).Value = $indexVals[$r]
}

$excel.CutCopyMode = $false
Write-Output "done"
